# Script to apply the target edit to the Turkey Super Lig 2023-2024 workbook.
#
# The edit consists of:
#   1. Swapping the match-data columns (B:V) between 15 pairs of adjacent
#      rows (this reorders some fixtures while keeping the running "Indice"
#      counter in column A untouched).
#   2. Appending three brand-new fixture rows (118-120) at the bottom of
#      the sheet, including proper styling on the Indice (A) and
#      data_partida (E) columns to match the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: swap columns B:V between each pair of rows listed below.
# ---------------------------------------------------------------------
$swapPairs = @(
    @(18, 19),
    @(29, 30),
    @(31, 32),
    @(33, 34),
    @(42, 43),
    @(44, 45),
    @(46, 47),
    @(49, 50),
    @(55, 56),
    @(63, 64),
    @(71, 72),
    @(82, 83),
    @(94, 95),
    @(97, 99),
    @(115, 116)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:V$r1")
    $range2 = $ws.Range("B$r2`:V$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value = $vals2
    $range2.Value = $vals1
}

# ---------------------------------------------------------------------
# Step 2: append the three new rows with fixture data.
# ---------------------------------------------------------------------
$newRows = @(
    @{ RowNum=118; Indice=117; DataPartida=45242.47916666666; Home="Pendikspor"; HomeGols=1; Away="Samsunspor"; AwayGols=0;
       HomeOpenOdds=2.71; HomeOpenDate="05/11/2023 15:44"; HomeCloseOdds=2.93; HomeCloseDate="12/11/2023 11:26";
       DrawOpenOdds=3.42; DrawOpenDate="05/11/2023 15:44"; DrawCloseOdds=3.51; DrawCloseDate="12/11/2023 11:26";
       AwayOpenOdds=2.71; AwayOpenDate="05/11/2023 15:44"; AwayCloseOdds=2.49; AwayCloseDate="12/11/2023 11:26";
       Url="https://www.betexplorer.com/football/turkey/super-lig/pendikspor-samsunspor/Kp3yw3dC/" },
    @{ RowNum=119; Indice=118; DataPartida=45242.58333333334; Home="Besiktas"; HomeGols=1; Away="Basaksehir"; AwayGols=0;
       HomeOpenOdds=1.55; HomeOpenDate="05/11/2023 17:12"; HomeCloseOdds=1.74; HomeCloseDate="12/11/2023 13:56";
       DrawOpenOdds=4.44; DrawOpenDate="05/11/2023 17:12"; DrawCloseOdds=3.84; DrawCloseDate="12/11/2023 13:59";
       AwayOpenOdds=5.79; AwayOpenDate="05/11/2023 17:12"; AwayCloseOdds=5.16; AwayCloseDate="12/11/2023 13:59";
       Url="https://www.betexplorer.com/football/turkey/super-lig/besiktas-basaksehir/vkz2gOtg/" },
    @{ RowNum=120; Indice=119; DataPartida=45242.58333333334; Home="Rizespor"; HomeGols=1; Away="Istanbulspor AS"; AwayGols=0;
       HomeOpenOdds=1.8; HomeOpenDate="06/11/2023 18:12"; HomeCloseOdds=1.69; HomeCloseDate="12/11/2023 13:53";
       DrawOpenOdds=3.94; DrawOpenDate="06/11/2023 18:12"; DrawCloseOdds=4.01; DrawCloseDate="12/11/2023 13:53";
       AwayOpenOdds=4.32; AwayOpenDate="06/11/2023 18:12"; AwayCloseOdds=5.3; AwayCloseDate="12/11/2023 13:52";
       Url="https://www.betexplorer.com/football/turkey/super-lig/rizespor-istanbulspor-as/04OAi2B5/" }
)

foreach ($row in $newRows) {
    $RowNum = $row.RowNum

    # Copy the number formatting used for the Indice and data_partida
    # columns from the last existing data row (117) so the new rows
    # match the look of the rest of the table.
    $ws.Range("A117").Copy() | Out-Null
    $ws.Range("A$RowNum").PasteSpecial(-4122) | Out-Null

    $ws.Range("E117").Copy() | Out-Null
    $ws.Range("E$RowNum").PasteSpecial(-4122) | Out-Null

    $ws.Range("A$RowNum").Value = $row.Indice
    $ws.Range("B$RowNum").Value = "turkey"
    $ws.Range("C$RowNum").Value = "super-lig"
    $ws.Range("D$RowNum").Value = "2023-2024"
    $ws.Range("E$RowNum").Value = $row.DataPartida
    $ws.Range("F$RowNum").Value = $row.Home
    $ws.Range("G$RowNum").Value = $row.HomeGols
    $ws.Range("H$RowNum").Value = $row.Away
    $ws.Range("I$RowNum").Value = $row.AwayGols
    $ws.Range("J$RowNum").Value = $row.HomeOpenOdds
    $ws.Range("K$RowNum").Value = $row.HomeOpenDate
    $ws.Range("L$RowNum").Value = $row.HomeCloseOdds
    $ws.Range("M$RowNum").Value = $row.HomeCloseDate
    $ws.Range("N$RowNum").Value = $row.DrawOpenOdds
    $ws.Range("O$RowNum").Value = $row.DrawOpenDate
    $ws.Range("P$RowNum").Value = $row.DrawCloseOdds
    $ws.Range("Q$RowNum").Value = $row.DrawCloseDate
    $ws.Range("R$RowNum").Value = $row.AwayOpenOdds
    $ws.Range("S$RowNum").Value = $row.AwayOpenDate
    $ws.Range("T$RowNum").Value = $row.AwayCloseOdds
    $ws.Range("U$RowNum").Value = $row.AwayCloseDate
    $ws.Range("V$RowNum").Value = $row.Url
}
